$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 5 (applied first, bottom-most edit, so earlier Find()-based lookups
# higher up the document are unaffected by this large restructure):
# Replace the block from "Click on connect & ok & yes..." through the
# "M.Madhuprasad..." signature paragraph with the corrected / tidied text,
# collapsing the run of blank paragraphs down to a single one and removing
# the stray lastRenderedPageBreak.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("Click on connect", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPara = $r1.Paragraphs(1)

$r2 = $d.Content
$r2.Find.Execute("M.Madhuprasad", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPara = $r2.Paragraphs(1)

$combined = $d.Range($startPara.Range.Start, $endPara.Range.End)
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Click on connect &amp; ok &amp; yes. And wait for a </w:t></w:r><w:r><w:t>while and</w:t></w:r><w:r><w:t xml:space="preserve"> verify.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Thanks &amp; Regards,</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>M.Madhuprasad</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$combined.InsertXML($xml5) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: drop the stray <w:lastRenderedPageBreak/> in front of
# "Now click on Launch instance."
# ---------------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Now click on Launch instance.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$launchPara = $r3.Paragraphs(1)
$launchPara.Range.Text = "Now click on Launch instance."

# ---------------------------------------------------------------------------
# Change 3: the "Answers:" paragraph right after "Login to VM as
# administrator." is split into an empty paragraph (keeping the ind+bold
# paragraph formatting) followed by a new paragraph carrying the bold
# "Answers:" run with a lastRenderedPageBreak in front of it.
# ---------------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Login to VM as administrator.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$loginVmPara = $r4.Paragraphs(1)
$answers2Para = $d.Paragraphs($loginVmPara.Index + 1)
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Answers:</w:t></w:r></w:p>
'@
$answers2Para.Range.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Change 2: insert a blank "ListParagraph"-styled paragraph and a bold
# "Topic 2: compute:" heading right after "Incase user is attached..." and
# before "Assignment 2:".
# ---------------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("Incase user is attached", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$incasePara = $r5.Paragraphs(1)
$idx = $incasePara.Index
$incasePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($idx + 1)
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Topic 2: compute:</w:t></w:r></w:p>
'@
$newPara.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Change 1: bold the first "Answers:" paragraph (both paragraph mark and
# the run) right after "Login to AWS console using TRHOL".
# ---------------------------------------------------------------------------
$r6 = $d.Content
$r6.Find.Execute("Login to AWS console using TRHOL", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$loginAwsPara = $r6.Paragraphs(1)
$answers1Para = $d.Paragraphs($loginAwsPara.Index + 1)
$answers1Para.Range.Font.Bold = 1

Write-Host "All changes applied"
